$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "Trabalho de pesquisa de referencias"
#           -> "Trabalho de pesquisa de requisitos"
# -----------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$r2.Find.Execute("referencias", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "requisitos", 2)

# -----------------------------------------------------------------
# Change 2: insert two new red paragraphs right before the
# paragraph that holds the _GoBack bookmark (the second-to-last
# paragraph; the very last paragraph is the trailing empty one).
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($count - 1)

# Insert a new (empty) paragraph mark right before the bookmark
# paragraph. After this call, $bookmarkPara refers to the newly
# created (now-first) empty paragraph, and the original bookmark
# paragraph has shifted to be its Next() sibling.
$bookmarkPara.Range.InsertParagraphBefore()
$newPara1 = $bookmarkPara
$realBookmarkPara = $newPara1.Next()

# New paragraph #1
$newPara1.Range.Text = "- cadastramento automático de placa, modelo, marca, tamanho e cor de veiculo por reconhecimento ótico"
$newPara1.Range.Font.Color = 255

# New paragraph #2's content is prepended into the bookmark
# paragraph itself, right before the bookmark start/end.
$insertPos = $realBookmarkPara.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore("- Cadastramento de politica de pontuação para descontos")

# Colour the whole bookmark paragraph (text + paragraph mark) red
$realBookmarkPara.Range.Font.Color = 255
